$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16 data
$ws.Range("A16").Value = 41184
$ws.Range("A16").NumberFormat = "ddd\ dd/mm/yyyy"

$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = "Doxygen set up, setEvent completed and first related, new test case tc05 succeeds"

# Update selection to match diff
$ws.Range("C16").Select()
